$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = "https://jobs.lever.co/StubHub/098afbcb-6412-4646-903f-df11e51f7b33"
$ws.Range("B2").Value = "location matched"

# Add new rows 3-7 with default search sites
$ws.Range("A3").Value = "https://jobs.boeing.com/category/engineering-software-jobs/185/2649/1"
$ws.Range("B3").Value = "Posting Website"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "

$ws.Range("A4").Value = "https://efds.fa.em5.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/requisitions"
$ws.Range("B4").Value = "Posting Website"
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "

$ws.Range("A5").Value = "https://nvidia.wd5.myworkdayjobs.com/en-US/NVIDIAExternalCareerSite/job/US-CA-Santa-Clara/Senior-Software-Engineer---HPC_JR1983439?locationHierarchy1=2fcb99c455831013ea52fb338f2932d8"
$ws.Range("B5").Value = "location matched"
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "

$ws.Range("A6").Value = "https://globalcareers-cotiviti.icims.com/jobs/12369/software-engineer/job?mobile=false&width=1100&height=500&bga=true&needsRedirect=false&jan1offset=330&jun1offset=330"
$ws.Range("B6").Value = "Posting Website"
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "

$ws.Range("A7").Value = "https://cooperative-sea-1e8.notion.site/d52db741a85748aead2235e7376f1974?v=902e7ec4bf28466f86dd89c8e9084427&pvs=74"
$ws.Range("B7").Value = "Posting Website"
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = " "
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = " "
